# CHRONICLER meeting-minutes doc: refresh date/time, resize the three
# inline chart images, and swap in the new opening/closing remarks.
$d = $word.ActiveDocument

# -- date / time stamps -----------------------------------------------
$d.Content.Find.Execute("날짜 : 2022-08-13", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "날짜 : 2022-08-14", 2) | Out-Null

$d.Content.Find.Execute("시간 : 17:03:03", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "시간 : 13:40:25", 2) | Out-Null

# -- speaker dialogue ---------------------------------------------------
$d.Content.Find.Execute("안녕하세요 춘자넷 및 아 반갑습니다", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "오늘 회의 시작하겠습니다", 2) | Out-Null

$d.Content.Find.Execute(" 정말 기쁩니다", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " 여러분들은 모두 최고예요", 2) | Out-Null

# -- inline chart image extents (EMU -> points, 12700 EMU per point) ----
# The three chart InlineShapes (wdInlineShapeChart = 12) appear, in document
# order, right after the title photo: docPr 100003 (word-frequency bar
# chart), 100004 (speaker pie chart) and 100005 (sentiment pie chart).
# Match on original size so this is robust to any reordering.
$chartShapes = @()
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    if ($shape.Type -eq 12) {
        $chartShapes += $shape
    }
}

# docPr id=100003 (word-frequency bar chart): 5486400 x 3200400 EMU
$shp1 = $chartShapes[0]
$shp1.Width  = 5969000 / 12700
$shp1.Height = 4114800 / 12700

# docPr id=100004 (speaker pie chart): 2730500 x 2794000 EMU
$shp2 = $chartShapes[1]
$shp2.Width  = 2984500 / 12700
$shp2.Height = 4114800 / 12700

# docPr id=100005 (sentiment pie chart): 2730500 x 2794000 EMU
$shp3 = $chartShapes[2]
$shp3.Width  = 2984500 / 12700
$shp3.Height = 4114800 / 12700
